$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds numeric-looking values stored as text; force text format
# so COM assignment does not coerce them into numbers (losing formatting
# like trailing/leading zeros, e.g. "262.50", "0.7000").
$priceCells = @(2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,18,19,21,23,27,40,41,42,43,44,45,46,47,48,49,50,51)
foreach ($r in $priceCells) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "262.23"

# Row 3
$ws.Range("D3").Value = "22.76"

# Row 4
$ws.Range("D4").Value = "6.209"

# Row 5
$ws.Range("D5").Value = "0.06114"

# Row 6
$ws.Range("D6").Value = "3.512"

# Row 7
$ws.Range("D7").Value = "6.705"

# Row 8
$ws.Range("D8").Value = "1.360"

# Row 9
$ws.Range("D9").Value = "0.7993"

# Row 10
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "0.01331"
$ws.Range("E10").Value = "9OneONE"

# Row 11
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "0.1569"
$ws.Range("E11").Value = "10WazirXWRX"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "0.08133"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"

# Row 13
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "0.03307"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"

# Row 14
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "0.03153"
$ws.Range("E14").Value = "13BitrueCoinBTR"

# Row 15
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "0.09271"
$ws.Range("E15").Value = "14BitMartTokenBMX"

# Row 16
$ws.Range("B16").Value = "MCDex"
$ws.Range("C16").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D16").Value = "3.903"
$ws.Range("E16").Value = "15MCDexMCB"

# Row 17
$ws.Range("B17").Value = "BitForexToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D17").Value = "0.001706"
$ws.Range("E17").Value = "16BitForexTokenBF"

# Row 18
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "0.04826"
$ws.Range("E18").Value = "17CoinExTokenCET"

# Row 19
$ws.Range("D19").Value = "0.006230"

# Row 21
$ws.Range("D21").Value = "0.003375"

# Row 23
$ws.Range("D23").Value = "3.693"

# Row 27
$ws.Range("D27").Value = "0.0004977"

# Row 40
$ws.Range("D40").Value = "0.04593"

# Row 41
$ws.Range("D41").Value = "0.007235"

# Row 42
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "0.003903"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "0.1119"
$ws.Range("E43").Value = "42BKEXTokenBKK"

# Row 44
$ws.Range("D44").Value = "0.01021"

# Row 45
$ws.Range("D45").Value = "0.002972"

# Row 46
$ws.Range("D46").Value = "0.00005955"

# Row 47
$ws.Range("D47").Value = "0.00000000751"

# Row 48
$ws.Range("D48").Value = "0.7005"

# Row 49
$ws.Range("D49").Value = "0.05048"
$ws.Range("E49").Value = "48BOLOBOLOWorstin24h"

# Row 50
$ws.Range("D50").Value = "0.00002101"

# Row 51
$ws.Range("D51").Value = "0.01011"
